$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9:125 down to 10:126
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with the new data record
$ws.Cells.Item(9, 1).Value = 4
$ws.Cells.Item(9, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(9, 3).Value = "Los Lagos"
$ws.Cells.Item(9, 4).Value = 44530
$ws.Cells.Item(9, 5).Value = 10
$ws.Cells.Item(9, 6).Value = 100112009
$ws.Cells.Item(9, 7).Value = "Acelga"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 250
$ws.Cells.Item(9, 11).Value = 3500
$ws.Cells.Item(9, 12).Value = 3500
$ws.Cells.Item(9, 13).Value = 3500
$ws.Cells.Item(9, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(9, 15).Value = "Región del Maule"
$ws.Cells.Item(9, 16).Value = 875
$ws.Cells.Item(9, 17).Value = 4
$ws.Cells.Item(9, 18).Value = "Hortaliza"
